$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordered list of "objet" values for rows 2-11 (header is row 1)
$objets = @(
    'Objet : Audit comptable et financier de l''ANEAQ pour l''ex 2023',
    'Objet : Acquisition des disjoncteurs 72,5 kV pour la Division Exploitation Transport Marrakech',
    'Objet : Peinture des lignes HTB du réseau de Division Exploitation Transport Agadir Lot n° 1 : Peinture des lignes 60 & 225 kV Lot n° 2 : Peinture des lignes 400 kV',
    'Objet : Exécution des travaux de régénération dans les périmètres relevant de la Direction Régionale de l’Agence Nationale des Eaux et Forêts de Fès Meknès repartis sur neuf lots ...',
    'Objet : Exécution des travaux de reboisement au profit de la Direction Régionale de l’Agence Nationale des Eaux et Forêts de Fès-Meknès répartis sur deux lots',
    'Objet : Exécution des travaux de sylviculture et conduite de peuplements forestiers dans les forêts relevant de la Direction Régionale de l’Agence Nationale des Eaux et Forêts Fès-Meknès, sur une superficie totale de 100 ha, la ...',
    'Objet : Entretien et réparation de matériels électriques aux stations piscicoles relevant du Centre National d''Hydrobiologie et de Pisciculture d’Azrou (lot unique).',
    'Objet : Achat et fourniture de matériels de campement destinés aux Parcs Nationaux relevant de la Direction Régionale de l’Agence Nationale des Eaux et Forêts de Fès Meknès répartis en deux lots ...',
    'Objet : Fourniture de la machine (AQUARIUS-TRUXOR DM 4700 B ou similaire) pour l’aménagement des milieux aquatiques relevant du Centre National d''Hydrobiologie et de Pisciculture (lot unique). ...',
    'Objet : Fourniture et pose de panneaux de signalisation au Parc National d’Ifrane relevant de la Direction Régionale de l’Agence Nationale des Eaux et Forêts de Fès-Meknès (Lot unique). ...'
)

for ($i = 0; $i -lt $objets.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $objets[$i]
    $ws.Cells.Item($row, 2).Value = "N/A"
}

